$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '44.023.98'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.358.99'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.84%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +2.64%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '235.27'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +1.26%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '72.67'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +10.99%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +21.05%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0989'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.89%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '28.34'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +6.28%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.715.04'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +2.22%  '
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +2.04%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '16.86'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +10.02%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.65'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +8.43%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.883'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +5.49%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.368.24'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +2.06%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '43.907.09'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +2.89%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '76.38'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +3.71%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.32'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +2.23%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '251.72'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.88%  '
$ws.Range('B23').NumberFormat = "@"
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').NumberFormat = "@"
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('B24').NumberFormat = "@"
$ws.Range('B24').Value = 'WEMIXToken'
$ws.Range('C24').NumberFormat = "@"
$ws.Range('C24').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.79'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.90%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '10.39'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +5.65%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.25'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.56%  '
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +1.49%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '173.19'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.57'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +8.83%  '
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +0.39%  '
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +4.50%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.20'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +4.49%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0711'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +3.72%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.17'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +3.86%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.74'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +2.43%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.44'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +2.25%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.43'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -1.47%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0269'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +6.76%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '19.35'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +11.24%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.93'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -1.35%  '
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +4.53%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0975'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +2.25%  '
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +1.76%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.45'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '97.89'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -0.72%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.182'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +12.74%  '
$ws.Range('B49').NumberFormat = "@"
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').NumberFormat = "@"
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.31'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('B50').NumberFormat = "@"
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').NumberFormat = "@"
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.438.14'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.586.39'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +2.07%  '
